# fix: superfluous next directive breaking the duplicate behavior
#
# The showcase table has one "record" cell per table row/column, each
# ending with a literal "{:next}" text run immediately followed by the
# real { NEXT } merge field. That literal run is superfluous everywhere,
# but removing it from the LAST record breaks the duplicate/merge
# behavior, so only the final occurrence is deleted.
#
# Find/Replace in this host always matches against the whole document
# content regardless of which Range invoked it, so instead we locate the
# exact paragraph whose entire text is "{:next}" (last one wins, which is
# the last paragraph in document order) and delete just that literal
# run's character span via an explicit Document.Range(start, end).

$d = $word.ActiveDocument

$needle = "{:next}"
$target = $null

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    if ($text -like ($needle + "*")) {
        # Keep overwriting so $target ends up the LAST matching paragraph.
        $target = $para
    }
}

if ($target -eq $null) {
    Write-Output "No '{:next}' paragraph found; nothing to do."
} else {
    $startPos = $target.Range.Start
    $endPos = $startPos + $needle.Length

    $victim = $d.Range($startPos, $endPos)
    if ($victim.Text -eq $needle) {
        $victim.Delete()
        Write-Output "Removed superfluous '{:next}' run."
    } else {
        Write-Output ("Unexpected text at target range: [" + $victim.Text + "]")
    }
}
